$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.308.17'
$ws.Range("E2").Value = '  +0.03%  '

$ws.Range("D3").Value = '2.006.48'
$ws.Range("E3").Value = '  -1.89%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").Value = "'251.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.66%  '

$ws.Range("D6").Value = "'0.638"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.29%  '

$ws.Range("D7").Value = "'61.38"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +10.50%  '

$ws.Range("E8").Value = '  -0.05%  '

$ws.Range("D9").Value = "'0.369"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.53%  '

$ws.Range("D10").Value = "'58.14"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.27%  '

$ws.Range("D11").Value = "'0.0740"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.33%  '

$ws.Range("E12").Value = '  -1.80%  '

$ws.Range("D13").Value = "'0.896"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.15%  '

$ws.Range("D14").Value = "'14.81"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.04%  '

$ws.Range("D15").Value = '2.305.34'
$ws.Range("E15").Value = '  -1.73%  '

$ws.Range("B16").Value = 'Polkadot'
$ws.Range("C16").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D16").Value = "'5.42"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.08%  '

$ws.Range("B17").Value = 'Avalanche'
$ws.Range("C17").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D17").Value = "'20.10"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +15.28%  '

$ws.Range("D18").Value = '2.015.99'
$ws.Range("E18").Value = '  -1.48%  '

$ws.Range("D19").Value = '36.314.15'
$ws.Range("E19").Value = '  +0.15%  '

$ws.Range("D20").Value = "'71.88"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.59%  '

$ws.Range("D21").Value = '0.0₃0859'
$ws.Range("E21").Value = '  +0.42%  '

$ws.Range("D22").Value = "'5.24"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.59%  '

$ws.Range("D23").Value = "'233.88"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.95%  '

$ws.Range("D24").Value = "'2.71"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +19.71%  '

$ws.Range("E25").Value = '  -0.18%  '

$ws.Range("D26").Value = "'2.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -1.37%  '

$ws.Range("D27").Value = "'9.53"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.73%  '

$ws.Range("D28").Value = "'163.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.47%  '

$ws.Range("D29").Value = "'19.54"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.03%  '

$ws.Range("E30").Value = '  -1.05%  '

$ws.Range("D31").Value = "'5.08"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.30%  '

$ws.Range("B32").Value = 'ImmutableX'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D32").Value = "'1.18"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.24%  '

$ws.Range("B33").Value = 'Kaspa'
$ws.Range("C33").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D33").Value = "'0.108"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +19.05%  '

$ws.Range("D34").Value = "'4.55"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +4.12%  '

$ws.Range("D35").Value = "'0.0606"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.13%  '

$ws.Range("D36").Value = "'2.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +10.15%  '

$ws.Range("E37").Value = '  -0.04%  '

$ws.Range("E38").Value = '  -1.08%  '

$ws.Range("D39").Value = "'5.90"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +16.83%  '

$ws.Range("D40").Value = "'0.104"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +14.48%  '

$ws.Range("B41").Value = 'RenderToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D41").Value = "'2.85"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +25.03%  '

$ws.Range("B42").Value = 'TrustWalletToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D42").Value = "'1.22"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.03%  '

$ws.Range("B43").Value = 'HuobiToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D43").Value = "'2.95"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.45%  '

$ws.Range("B44").Value = 'FraxShare'
$ws.Range("C44").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D44").Value = "'8.10"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +8.17%  '

$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").Value = "'0.0215"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.30%  '

$ws.Range("E46").Value = '  +2.47%  '

$ws.Range("B47").Value = 'InjectiveProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D47").Value = "'16.91"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +7.43%  '

$ws.Range("D48").Value = '1.450.10'
$ws.Range("E48").Value = '  +3.52%  '

$ws.Range("D49").Value = "'94.35"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.39%  '

$ws.Range("D50").Value = "'2.93"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.98%  '

$ws.Range("D51").Value = "'46.93"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +1.33%  '
